# Inserts a new weekly price record (Zapallo italiano, Macroferia Regional de
# Talca) as a new row 367, shifting the existing rows 367:396 down to 368:397.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 367; Excel shifts rows 367:396 -> 368:397
# and the sheet dimension grows from A1:R396 to A1:R397 automatically.
$ws.Rows.Item(367).Insert()

# Populate the new row 367 with the new record's data.
$ws.Cells.Item(367, 1).Value = 5
$ws.Cells.Item(367, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(367, 3).Value = 'Maule'
$ws.Cells.Item(367, 4).Value = 44783
$ws.Cells.Item(367, 5).Value = 7
$ws.Cells.Item(367, 6).Value = 100112032
$ws.Cells.Item(367, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(367, 8).Value = 'Sin especificar'
$ws.Cells.Item(367, 9).Value = 'Primera'
$ws.Cells.Item(367, 10).Value = 300
$ws.Cells.Item(367, 11).Value = 20000
$ws.Cells.Item(367, 12).Value = 20000
$ws.Cells.Item(367, 13).Value = 20000
$ws.Cells.Item(367, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(367, 15).Value = 'Región del Maule'
$ws.Cells.Item(367, 16).Value = 400
$ws.Cells.Item(367, 17).Value = 50
$ws.Cells.Item(367, 18).Value = 'Hortaliza'
